$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item(4).Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item(2).Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
Write-Host "before rename:"
foreach ($ws in $wb.Worksheets) { Write-Host ($ws.Name + " | " + $ws.Index) }
$wb.Worksheets.Item(4).Name = "Building-List OLD"
$wb.Worksheets.Item(1).Name = "Building-List"
Write-Host "after rename:"
foreach ($ws in $wb.Worksheets) { Write-Host ($ws.Name + " | " + $ws.Index) }
